$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.362.98'
$ws.Range('E2').Value = '  +6.07%  '
$ws.Range('D3').Value = '2.365.60'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'110.07"
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').Value = "'310.30"
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.618"
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').Value = "'41.34"
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = "'8.48"
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = "'0.984"
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '2.722.12'
$ws.Range('E15').Value = '  +2.17%  '
$ws.Range('D16').Value = "'15.38"
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '2.359.03'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '45.315.61'
$ws.Range('E18').Value = '  +5.97%  '
$ws.Range('D19').Value = "'7.32"
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = "'13.73"
$ws.Range('E21').Value = '  +4.59%  '
$ws.Range('D22').Value = "'73.31"
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = "'3.45"
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = "'259.04"
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('D25').Value = "'2.30"
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').Value = "'11.16"
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('D28').Value = "'7.39"
$ws.Range('E28').Value = '  -3.61%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = "'0.0972"
$ws.Range('E30').Value = '  +11.31%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'38.19"
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'22.42"
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = "'170.50"
$ws.Range('D34').Value = "'2.92"
$ws.Range('E34').Value = '  +6.56%  '
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').Value = "'4.84"
$ws.Range('E36').Value = '  +4.17%  '
$ws.Range('D37').Value = "'0.115"
$ws.Range('E37').Value = '  +1.85%  '
$ws.Range('D38').Value = "'2.97"
$ws.Range('E38').Value = '  +4.66%  '
$ws.Range('D39').Value = "'3.94"
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = "'1.74"
$ws.Range('E41').Value = '  +8.90%  '
$ws.Range('D42').Value = "'99.03"
$ws.Range('E42').Value = '  -5.55%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').Value = "'69.97"
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('D45').Value = "'12.85"
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = "'82.57"
$ws.Range('E47').Value = '  +7.69%  '
$ws.Range('D48').Value = "'113.13"
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').Value = "'9.26"
$ws.Range('E49').Value = '  +4.47%  '
$ws.Range('D50').Value = "'5.49"
$ws.Range('E50').Value = '  +4.54%  '
$ws.Range('D51').Value = '1.667.88'
$ws.Range('E51').Value = '  -0.20%  '
